$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4484.9443
$ws.Range("I40").Value = 4825.385
$ws.Range("K40").Value = 4825.385
$ws.Range("M40").Value = -4650.385
$ws.Range("H93").Value = 100600.5
$ws.Range("J93").Value = 100600.5
$ws.Range("L93").Value = 100600.5
$ws.Range("N93").Value = -105592.5
$ws.Range("H113").Value = 4479.067
$ws.Range("I113").Value = 4366.778
$ws.Range("J113").Value = 4647.5
$ws.Range("K113").Value = 4366.778
$ws.Range("L113").Value = 4647.5
$ws.Range("M113").Value = -1112.778
$ws.Range("N113").Value = -11155.5
$ws.Range("H116").Value = 91578.11
$ws.Range("I116").Value = 115457.57
$ws.Range("K116").Value = 115457.57
$ws.Range("M116").Value = -112015.57
$ws.Range("H132").Value = 6247.396
$ws.Range("I132").Value = 6214.9565
$ws.Range("K132").Value = 18644.8695
$ws.Range("M132").Value = -16114.8695
$ws.Range("H138").Value = 3566.868
$ws.Range("I138").Value = 1592.8
$ws.Range("J138").Value = 4025.9534
$ws.Range("K138").Value = 4778.4
$ws.Range("L138").Value = 12077.8602
$ws.Range("M138").Value = 361.6000000000004
$ws.Range("N138").Value = -22357.8602
$ws.Range("H140").Value = 75000
$ws.Range("J140").Value = 75000
$ws.Range("L140").Value = 75000
$ws.Range("N140").Value = -85360

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4011.4167
$ws.Range("I61").Value = 3571
$ws.Range("K61").Value = 3571
$ws.Range("M61").Value = -3359
$ws.Range("H122").Value = 4546.75
$ws.Range("I122").Value = 3750
$ws.Range("K122").Value = 11250
$ws.Range("M122").Value = -8800
$ws.Range("H130").Value = 49095
$ws.Range("J130").Value = 49095
$ws.Range("L130").Value = 49095
$ws.Range("N130").Value = -59135
$ws.Range("H136").Value = 4011.4167
$ws.Range("I136").Value = 3571
$ws.Range("K136").Value = 10713
$ws.Range("M136").Value = -8163

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 19362.375
$ws.Range("J75").Value = 14985.714
$ws.Range("L75").Value = 14985.714
$ws.Range("N75").Value = -16857.714
$ws.Range("H78").Value = 19362.375
$ws.Range("J78").Value = 14985.714
$ws.Range("L78").Value = 44957.142
$ws.Range("N78").Value = -54317.142
$ws.Range("H99").Value = 39111.645
$ws.Range("I99").Value = 53086.8
$ws.Range("J99").Value = 4173.75
$ws.Range("K99").Value = 53086.8
$ws.Range("L99").Value = 4173.75
$ws.Range("M99").Value = -51588.8
$ws.Range("N99").Value = -7169.75
$ws.Range("H128").Value = 12283.8
$ws.Range("I128").Value = 12283.8
$ws.Range("K128").Value = 36851.39999999999
$ws.Range("M128").Value = -34361.39999999999
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()
$ws.Range("H134").Value = 14700.277
$ws.Range("I134").Value = 15287.8125
$ws.Range("K134").Value = 45863.4375
$ws.Range("M134").Value = -43328.4375

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2598.25
$ws.Range("I16").Value = 2431
$ws.Range("K16").Value = 2431
$ws.Range("M16").Value = -2144
$ws.Range("H99").Value = 4789.15
$ws.Range("I99").Value = 6999.6665
$ws.Range("J99").Value = 4399.0586
$ws.Range("K99").Value = 6999.6665
$ws.Range("L99").Value = 4399.0586
$ws.Range("M99").Value = -5501.6665
$ws.Range("N99").Value = -7395.0586
$ws.Range("H113").Value = 2598.25
$ws.Range("I113").Value = 2431
$ws.Range("K113").Value = 2431
$ws.Range("M113").Value = -261
$ws.Range("H126").Value = 4789.15
$ws.Range("I126").Value = 6999.6665
$ws.Range("J126").Value = 4399.0586
$ws.Range("K126").Value = 20998.9995
$ws.Range("L126").Value = 13197.1758
$ws.Range("M126").Value = -18528.9995
$ws.Range("N126").Value = -18137.1758

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 6016.5264
$ws.Range("J93").Value = 6128.5557
$ws.Range("L93").Value = 18385.6671
$ws.Range("N93").Value = -22129.6671
$ws.Range("H131").Value = 2002.1923
$ws.Range("J131").Value = 3684.6365
$ws.Range("L131").Value = 11053.9095
$ws.Range("N131").Value = -21133.9095
$ws.Range("H137").Value = 4604.077
$ws.Range("J137").Value = 6622.75
$ws.Range("L137").Value = 19868.25
$ws.Range("N137").Value = -30068.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3807.7144
$ws.Range("I113").Value = 2914.125
$ws.Range("J113").Value = 4999.1665
$ws.Range("K113").Value = 2914.125
$ws.Range("L113").Value = 4999.1665
$ws.Range("M113").Value = -744.125
$ws.Range("N113").Value = -9339.166499999999
$ws.Range("H126").Value = 3900
$ws.Range("J126").Value = 4311.1113
$ws.Range("L126").Value = 12933.3339
$ws.Range("N126").Value = -17873.3339
$ws.Range("H132").Value = 6531.5
$ws.Range("I132").Value = 6664.737
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 19994.211
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -17464.211
$ws.Range("N132").Value = -17060

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 6249.5
$ws.Range("J5").Value = 6249.5
$ws.Range("L5").Value = 6249.5
$ws.Range("N5").Value = -6475.5
$ws.Range("H7").Value = 3454.9092
$ws.Range("I7").Value = 3444.889
$ws.Range("J7").Value = 3500
$ws.Range("K7").Value = 3444.889
$ws.Range("L7").Value = 3500
$ws.Range("M7").Value = -3332.889
$ws.Range("N7").Value = -3724
$ws.Range("H31").Value = 1516.4615
$ws.Range("I31").Value = 2541.4285
$ws.Range("J31").Value = 320.66666
$ws.Range("K31").Value = 2541.4285
$ws.Range("L31").Value = 320.66666
$ws.Range("M31").Value = -2293.4285
$ws.Range("N31").Value = -816.66666
$ws.Range("H35").Value = 2627.625
$ws.Range("J35").Value = 13999.333
$ws.Range("L35").Value = 13999.333
$ws.Range("N35").Value = -14671.333
$ws.Range("H43").Value = 184166.33
$ws.Range("I43").Value = 5000
$ws.Range("K43").Value = 5000
$ws.Range("M43").Value = -4807
$ws.Range("H46").Value = 3567.48
$ws.Range("J46").Value = 4049.5
$ws.Range("L46").Value = 4049.5
$ws.Range("N46").Value = -4425.5
$ws.Range("H98").Value = 100355
$ws.Range("J98").Value = 100355
$ws.Range("L98").Value = 100355
$ws.Range("N98").Value = -106345
$ws.Range("H100").Value = 13555.223
$ws.Range("I100").Value = 2199.8
$ws.Range("K100").Value = 2199.8
$ws.Range("M100").Value = -1658.8
$ws.Range("H126").Value = 3454.9092
$ws.Range("I126").Value = 3444.889
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 10334.667
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -7864.667000000001
$ws.Range("N126").Value = -15440
$ws.Range("H130").Value = 96571.75
$ws.Range("J130").Value = 96571.75
$ws.Range("L130").Value = 96571.75
$ws.Range("N130").Value = -106611.75
$ws.Range("H132").Value = 29663.3
$ws.Range("I132").Value = 34704.125
$ws.Range("K132").Value = 104112.375
$ws.Range("M132").Value = -101582.375
$ws.Range("H136").Value = 39750.77
$ws.Range("I136").Value = 3929.9
$ws.Range("K136").Value = 11789.7
$ws.Range("M136").Value = -9239.700000000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 7000
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H36").Value = 7000
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H92").Value = 112511
$ws.Range("J92").Value = 112511
$ws.Range("L92").Value = 112511
$ws.Range("N92").Value = -117503
$ws.Range("H95").Value = 60285.832
$ws.Range("J95").Value = 60285.832
$ws.Range("L95").Value = 60285.832
$ws.Range("N95").Value = -65777.83199999999
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H130").Value = 85214.5
$ws.Range("J130").Value = 85214.5
$ws.Range("L130").Value = 85214.5
$ws.Range("N130").Value = -95254.5
